$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# --- Header row: add translation/locale setting columns ---
$ws.Range("D1").Value = "display.title.text.pt"
$ws.Range("E1").Value = "display.title.text.sw"
$ws.Range("F1").Value = "display.locale.text"
$ws.Range("G1").Value = "display.locale.text.pt"
$ws.Range("H1").Value = "display.locale.text.sw"

# --- Row 5 (survey / display title): mirror title into pt/sw columns ---
$ws.Range("D5").Value = $ws.Range("C5").Value2
$ws.Range("E5").Value = $ws.Range("C5").Value2

# --- New locale rows: default (English), pt (Português), sw (Kiswahili) ---
$ws.Range("A7").Value = "default"
$ws.Range("F7").Value = "English"
$ws.Range("G7").Value = "English"
$ws.Range("H7").Value = "English"

$ws.Range("A8").Value = "pt"
$ws.Range("F8").Value = "Português"
$ws.Range("G8").Value = "Português"
$ws.Range("H8").Value = "Português"

$ws.Range("A9").Value = "sw"
$ws.Range("F9").Value = "Kiswahili"
$ws.Range("G9").Value = "Kiswahili"
$ws.Range("H9").Value = "Kiswahili"

# --- Make "settings" the active sheet / selected cell (was "model") ---
$ws.Range("A1").Select() | Out-Null
$ws.Activate() | Out-Null
